# Rename the sheet: "Sheet1" -> "Sheet"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Sheet"

# The sheet used to hold 5 rows (A2:A6) of the same repeated placeholder
# string. Replace that with a small 2-row x 3-col table of real data
# (pycard output: a couple of hashes, a name, and a count per row).
$ws.Range("A2:A6").ClearContents()

$ws.Range("A2").Value = "CFD893A460"
$ws.Range("B2").Value = "gary tsai"
$ws.Range("C2").Value = 6
$ws.Range("A3").Value = "8A9AB340"
$ws.Range("B3").Value = "CBS"
$ws.Range("C3").Value = 1

# Match the selection left behind by the edit.
$ws.Range("A2:C3").Select() | Out-Null
